$d = $word.ActiveDocument

# The "Predicting fingers" section ends with problem 1 (a blank, left-indented
# paragraph is used as a content placeholder for subsequent items). The
# paragraph immediately following "a) The problem is finding out what
# finger ..." is a plain blank paragraph, and the one after that is the
# first of three blank, left-indented (360 twips) paragraphs. That first
# blank indented paragraph is the one that gets filled in with the new
# "over all goal" content for problem 3 below.
$anchor = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -match "The problem is finding out what finger") {
        $anchor = $i
        break
    }
}

if ($anchor -eq $null) {
    throw "edit.ps1: could not locate the 'Predicting fingers' problem 1 paragraph"
}

$target = $d.Paragraphs.Item($anchor + 2)

if ($target.Range.Text.Trim() -ne "" -or $target.Format.LeftIndent -eq 0) {
    throw "edit.ps1: target paragraph is not the expected blank, indented placeholder"
}

# Replace that blank paragraph with four new paragraphs:
#   - two numbered ("a)") list-paragraph bullets (sharing the existing
#     numId=2 list used elsewhere in the document) giving the "over all
#     goal" text for the fingers problem
#   - a plain "2.  " heading paragraph
#   - a plain "a)" paragraph, ready for the next answer
$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">You have to get the number of count for all fingers counted </w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Over all goal is to find a </w:t></w:r><w:r><w:t xml:space="preserve">problem </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>sentece</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> to find the finger where all counts asked land on </w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:ind w:left="360"/></w:pPr><w:r><w:t xml:space="preserve">2.  </w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:ind w:left="360"/></w:pPr><w:proofErr w:type="gramStart"/><w:r><w:t>a</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>)</w:t></w:r></w:p>
'@

$target.Range.InsertXML($xml)

Write-Host "Inserted over-all-goal paragraphs for problem 3 (fingers) after paragraph $anchor"
